# Reconfiguration Experiments workbook - processing of evaluation results.
#
# The three sheets (Experiment_1, Experiement_2, Experiment_3) each contain a
# "percentage" mirror of the main metrics table (row 12 mirrors the row 6
# headers, row 26 mirrors the row 20/23 headers). Those mirrored header cells
# are updated to explicit "(%)" labelled headers so they read correctly as a
# percentage breakdown rather than re-using the raw absolute-value labels.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Experiment_1", "Experiement_2", "Experiment_3")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 12: percentage-table header mirror of row 6.
    $ws.Range("B12").Value = "Unity code (%)"
    $ws.Range("C12").Value = "Services Code (%)"
    $ws.Range("D12").Value = "DT code (%)"
    $ws.Range("E12").Value = "Configuration File (%)"
    $ws.Range("F12").Value = "Percentage Code Reused (%)"

    # Row 26: percentage-table header mirror of row 20/23.
    $ws.Range("B26").Value = "Source Code (%)"
    $ws.Range("C26").Value = "Configuration File (%)"
    $ws.Range("D26").Value = "Percentage Code Reused (%)"
}

# Update the active sheet/selection bookkeeping to match the latest review
# pass: Experiment_1 becomes the active tab (cursor left on F12, the cell
# just updated), Experiement_2's cursor is left on D26, and Experiment_3
# (previously the active tab) is no longer active, with its cursor left on
# B23.
$ws2 = $wb.Worksheets.Item("Experiement_2")
$ws2.Select() | Out-Null
$ws2.Range("D26").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("Experiment_3")
$ws3.Select() | Out-Null
$ws3.Range("B23").Select() | Out-Null

$ws1 = $wb.Worksheets.Item("Experiment_1")
$ws1.Select() | Out-Null
$ws1.Range("F12").Select() | Out-Null
